$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet/tab title text (the workbook.xml <sheet name="..."> entry)
$ws.Name = "BOM_Special_Special remake_1_20"

# 2. Row 4 (No.3, 100nF capacitor): quantity 1 -> 2, designator C2 -> C2,C3
$ws.Range("B4").Value = 2
$ws.Range("D4").Value = "C2,C3"

# 3. Row 6 (No.5, CONN1 header): add missing Comment value
$ws.Range("C6").Value = "HDR-M_2.54_1x10P"

# 4. Remove the transistor row (No.11, Q1 / UMH3N) entirely; everything below shifts up one row
$ws.Rows.Item(12).Delete()

# 5. The resistor row that used to be No.12 is now in row 12 (labelled No.11):
#    quantity 4 -> 5, designator gains R10
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 5
$ws.Range("D12").Value = "R1,R2,R3,R10,RD"

# 6. Renumber the remaining rows (previously No.13..No.21, now shifted up to rows 13..21)
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18
$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20

# 7. Row 18 (now No.17) switches part from CH340K/UART to CH340X/UART1
$ws.Range("C18").Value = "CH340X"
$ws.Range("D18").Value = "UART1"
$ws.Range("E18").Value = "MSOP-10_L3.0-W3.0-P0.50-LS5.0-BL"
$ws.Range("G18").Value = "CH340X"
$ws.Range("I18").Value = "C3035748"
